$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, copying the formatting used by the existing
# header row (bold, centered, bordered) from the last header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
